$d = $word.ActiveDocument

# --- 1. Title paragraph: "Symbol Graphs" -> "4 " + "Symbol Graphs" as two runs ---
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="3519BA25" w14:textId="0533502A" w:rsidR="00FC4C8C" w:rsidRDefault="009302FF" w:rsidP="009302FF">' +
  '<w:pPr><w:pStyle w:val="IntenseQuote"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">4 </w:t></w:r>' +
  '<w:r><w:t>Symbol Graphs</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($titleXml) | Out-Null

# --- 2. Merge the three runs of the "Typical applications..." paragraph into one run ---
$d.Content.Find.Execute(
    "Typical applications involve processing graphs defined in files or on web pages, using strings, not integer indices, to define and refer to vertices.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Typical applications involve processing graphs defined in files or on web pages, using strings, not integer indices, to define and refer to vertices.",
    2) | Out-Null

# --- 3. Mark the three picture runs as NoProofing (adds <w:rPr><w:noProof/></w:rPr>) ---
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = 1
}
